# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = -0.113
$ws.Range("W2").Value = 1.153061224489796
$ws.Range("X2").Value = 0.06341144098668264
$ws.Range("Y2").Value = 1.089649783503113
$ws.Range("AA2").Value = 1.298850574712644
$ws.Range("AB2").Value = 0.06228762779337538
$ws.Range("AC2").Value = 1.236562946919268
$ws.Range("AD2").Value = 0.176
$ws.Range("AF2").Value = 0.176
$ws.Range("AG2").Value = 0.176
$ws.Range("AH2").Value = 0.123422159887798
$ws.Range("AI2").Value = -5.333333333333333
$ws.Range("AJ2").Value = 0.123422159887798
$ws.Range("AK2").Value = -5.333333333333333

$ws.Range("K3").Value = -0.113
$ws.Range("W3").Value = 1.153061224489796
$ws.Range("X3").Value = 0.06341144098668264
$ws.Range("Y3").Value = 1.089649783503113
$ws.Range("AA3").Value = 1.298850574712644
$ws.Range("AB3").Value = 0.06228762779337538
$ws.Range("AC3").Value = 1.236562946919268
$ws.Range("AD3").Value = 0.176
$ws.Range("AF3").Value = 0.176
$ws.Range("AG3").Value = 0.176
$ws.Range("AH3").Value = 0.123422159887798
$ws.Range("AI3").Value = -5.333333333333333
$ws.Range("AJ3").Value = 0.123422159887798
$ws.Range("AK3").Value = -5.333333333333333
